$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1786
$ws1.Range("F7").Value = 323
$ws1.Range("F8").Value = 508
$ws1.Range("F9").Value = 4645
$ws1.Range("F11").Value = 452
$ws1.Range("F13").Value = 987
$ws1.Range("F17").Value = 2987
$ws1.Range("F20").Value = 45
$ws1.Range("F22").Value = 36
$ws1.Range("F24").Value = 930
$ws1.Range("F27").Value = 2617
$ws1.Range("F28").Value = 1017
$ws1.Range("F29").Value = 2456
$ws1.Range("F31").Value = 1330
$ws1.Range("F32").Value = 608
$ws1.Range("F33").Value = 91
$ws1.Range("F34").Value = 884
$ws1.Range("F35").Value = 427
$ws1.Range("F36").Value = 1125
$ws1.Range("F37").Value = 924
$ws1.Range("F38").Value = 1180
$ws1.Range("F39").Value = 18
$ws1.Range("F40").Value = 860
$ws1.Range("F41").Value = 523
$ws1.Range("F42").Value = 362
$ws1.Range("F43").Value = 284
$ws1.Range("F44").Value = 3488

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 884

# Sheet "全部类型" (all types, aggregated)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 1786
$ws4.Range("F8").Value = 323
$ws4.Range("F9").Value = 508
$ws4.Range("F10").Value = 4645
$ws4.Range("F15").Value = 2987
$ws4.Range("F19").Value = 45
$ws4.Range("F21").Value = 884
$ws4.Range("F24").Value = 36
$ws4.Range("F26").Value = 930
$ws4.Range("F28").Value = 2617
$ws4.Range("F31").Value = 1017
$ws4.Range("F32").Value = 2456
$ws4.Range("F33").Value = 1330
$ws4.Range("F34").Value = 608
$ws4.Range("F36").Value = 91
$ws4.Range("F37").Value = 884
$ws4.Range("F38").Value = 1125
$ws4.Range("F39").Value = 924
$ws4.Range("F41").Value = 1180
$ws4.Range("F42").Value = 860
$ws4.Range("F43").Value = 523
$ws4.Range("F44").Value = 362
$ws4.Range("F48").Value = 284
$ws4.Range("F49").Value = 3488

$wb.Save()
